$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append after the existing data (rows 168-170)
$rows = @(
    @{ Row = 168; A = "2024-05-14"; B = "18:25:34"; C = "No pone tornillo"; D = "-"; E = "-"; F = "-"; G = "-" },
    @{ Row = 169; A = "2024-05-14"; B = "18:25:42"; C = "-"; D = "-"; E = "-"; F = "Fallo atornillador"; G = "-" },
    @{ Row = 170; A = "2024-05-14"; B = "18:31:21"; C = "-"; D = "Cámara no detecta Power CP"; E = "-"; F = "-"; G = "-" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Columns A (date) and B (time) must stay as plain text, not be
    # auto-converted to date/time serial numbers.
    $ws.Cells.Item($rowNum, 1).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 1).Value = $r.A

    $ws.Cells.Item($rowNum, 2).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 2).Value = $r.B

    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
}
